# Update the "urban"/"rural" breakdown row labels (A6:C6 and A7:C7) to the
# new, more formal wording, and move the active selection.
#
# 1.1.1.1a Level of extreme poverty sheet — "Kyrgyz Republic" row totals are
# followed by an urban/rural split. The short labels (шаар/город/urban and
# айыл/село/rural) are replaced with the longer forms used elsewhere in the
# platform (Шаар жерлери/Городские поселения/City and
# Айыл аймагы/Сельская местность/Village).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "Шаар жерлери"
$ws.Range("B6").Value = "Городские поселения"
$ws.Range("C6").Value = "City"

$ws.Range("A7").Value = "Айыл аймагы"
$ws.Range("B7").Value = "Сельская местность"
$ws.Range("C7").Value = "Village"

# Move / update the saved selection on the sheet (was M20, now B29) and
# drop the scrolled-right viewport (was topLeftCell = C1).
$ws.Range("B29").Select() | Out-Null
